$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.178490877151489
$ws.Range("B1").Value = 2.416741847991943
$ws.Range("D1").Value = 2.332604646682739
$ws.Range("E1").Value = 1.197142362594604
